# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text would otherwise be auto-converted
# to a real number by Excel; force them to keep Text format first.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D16', 'D18', 'D20', 'D21', 'D22', 'D25', 'D26', 'D27', 'D28', 'D29', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Cell value updates (matches the refreshed cryptos snapshot)
$ws.Range('D2').Value = '29.903.41'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.892.97'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '0.7732'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').Value = '243.92'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '25.71'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('D10').Value = '0.07236'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').Value = '0.08741'
$ws.Range('E11').Value = '  +7.95%  '
$ws.Range('D12').Value = '2.069.54'
$ws.Range('E12').Value = '  +7.32%  '
$ws.Range('D13').Value = '0.7717'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '5.407'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '94.28'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '6.220'
$ws.Range('E16').Value = '  +0.89%  '
$ws.Range('D17').Value = '30.100.45'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '13.94'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = '2.325.88'
$ws.Range('E19').Value = '  +8.69%  '
$ws.Range('D20').Value = '245.54'
$ws.Range('D21').Value = '0.000007868'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').Value = '8.212'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '0.1598'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').Value = '9.535'
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('D27').Value = '162.62'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D28').Value = '18.82'
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('D29').Value = '2.046'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').Value = '4.522'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').Value = '4.124'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '0.05468'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').Value = '1.250'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('D36').Value = '0.7533'
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('D37').Value = '1.005'
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').Value = '2.690'
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('D39').Value = '0.01982'
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').Value = '2.786'
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('D41').Value = '0.4514'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').Value = '74.07'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.101.49'
$ws.Range('E43').Value = '  -3.82%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '6.087'
$ws.Range('E44').Value = '  +3.82%  '
$ws.Range('D45').Value = '0.8551'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('D46').Value = '2.194.13'
$ws.Range('E46').Value = '  +7.27%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').Value = '102.96'
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('D49').Value = '1.886'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '7.622'
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').Value = '9.835'
$ws.Range('E51').Value = '  -1.60%  '
